$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

# The sheet repeats a 4-row "sprint summary" block (header + Written/Execution/
# Review rows). Add the next one ("Spint( 34) - Day 5- Test Case Summary") by
# cloning the most recent block (B41:C44) down to B47:C50 so the new block
# inherits the same cell styles / the merged header cell, then relabel it and
# fill in its numbers.
$ws.Range("B41:C44").Copy($ws.Range("B47:C50"))
$excel.CutCopyMode = $false

$ws.Range("B47").Value = "Spint( 34) - Day 5- Test Case Summary"
$ws.Range("C48").Value = 80
$ws.Range("C49").Value = 52
$ws.Range("C50").Value = 0

# Rows 41:44 render at 18pt (same font-driven height as the rest of the
# sprint blocks); give the new rows the same height.
$ws.Rows.Item(47).RowHeight = 18
$ws.Rows.Item(48).RowHeight = 18
$ws.Rows.Item(49).RowHeight = 18
$ws.Rows.Item(50).RowHeight = 18

# Scroll the view down to the newly added block and select it.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1
$ws.Range("H41").Select()
